# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds the "K" (strikeout) stat for each start, recomputed here
# to replace the old Strike# counts with the corrected K values.
$kValues = @{
    2 = 1
    3 = 1
    4 = 3
    5 = 3
    6 = 6
    7 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
